$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two now-unused columns (firing_degree_min, firing_degree_max);
# the remaining firing_degree column (J) shifts left into H.
$ws.Columns("H:I").Delete()

# Rewrite header row (B1:H1)
$ws.Cells.Item(1, 2).Value = "mean"
$ws.Cells.Item(1, 3).Value = "std"
$ws.Cells.Item(1, 4).Value = "y_mean"
$ws.Cells.Item(1, 5).Value = "y_std"
$ws.Cells.Item(1, 6).Value = "NumObservations"
$ws.Cells.Item(1, 7).Value = "tau"
$ws.Cells.Item(1, 8).Value = "firing_degree"

# Rewrite data rows 2..20 (A..H), replacing previous 16-row block with 19 rows
$data = @(
  @("0", "[[0.05311143]`n [0.1302186 ]]", "[[0.0849606 ]`n [0.09574303]]", [double]"-16.4835215262267", [double]"3.252131102432681", [double]"29", [double]"0.002097373084493301", [double]"8.531067180595293e-06"),
  @("1", "[[0.08228968]`n [0.15556817]]", "[[0.01763293]`n [0.07905498]]", [double]"-15.26825983559251", [double]"0.5728157217851149", [double]"138", [double]"1e-10", [double]"4.067501029582579e-13"),
  @("2", "[[0.13467593]`n [0.19719347]]", "[[0.01996375]`n [0.09478911]]", [double]"-13.17499095679213", [double]"0.5934112218553106", [double]"202", [double]"1e-10", [double]"4.067501029582579e-13"),
  @("3", "[[0.18632657]`n [0.24335731]]", "[[0.02015919]`n [0.09195827]]", [double]"-11.07896038543913", [double]"0.6101912906905872", [double]"314", [double]"1e-10", [double]"4.067501029582579e-13"),
  @("4", "[[0.24124181]`n [0.28499409]]", "[[0.01698617]`n [0.06268645]]", [double]"-8.891469242795642", [double]"0.5796707419312355", [double]"767", [double]"1.048508680493986e-10", [double]"4.26481013743556e-13"),
  @("5", "[[0.28659431]`n [0.32016954]]", "[[0.01852733]`n [0.06793296]]", [double]"-7.081710794256293", [double]"0.6193674930279515", [double]"585", [double]"0.01721672941677837", [double]"7.002906462879071e-05"),
  @("6", "[[0.34163201]`n [0.36370714]]", "[[0.01898801]`n [0.07252254]]", [double]"-4.881486813462899", [double]"0.609263568790357", [double]"419", [double]"123.572293815723", [double]"0.5026304323233344"),
  @("7", "[[0.39534978]`n [0.4024417 ]]", "[[0.01805944]`n [0.06282955]]", [double]"-2.755037955733696", [double]"0.6192359051777768", [double]"457", [double]"122.246803365469", [double]"0.4972389985522241"),
  @("8", "[[0.44815451]`n [0.44038909]]", "[[0.01800461]`n [0.06362721]]", [double]"-0.6676644234046379", [double]"0.6181938485618268", [double]"380", [double]"0.01278647172802157", [double]"5.200898691845625e-05"),
  @("9", "[[0.50108669]`n [0.47953589]]", "[[0.01742155]`n [0.06115823]]", [double]"1.432822925176326", [double]"0.6034796269726226", [double]"420", [double]"1e-10", [double]"4.067501029582579e-13"),
  @("10", "[[0.5509719 ]`n [0.51549336]]", "[[0.01756002]`n [0.06402713]]", [double]"3.410448595297656", [double]"0.5928423466914202", [double]"483", [double]"1e-10", [double]"4.067501029582579e-13"),
  @("11", "[[0.6044436 ]`n [0.55615225]]", "[[0.01883843]`n [0.0781522 ]]", [double]"5.537533146667895", [double]"0.5940084649370668", [double]"403", [double]"1e-10", [double]"4.067501029582579e-13"),
  @("12", "[[0.65830094]`n [0.59596906]]", "[[0.01951538]`n [0.0855098 ]]", [double]"7.676381887746297", [double]"0.5954232679831231", [double]"430", [double]"1e-10", [double]"4.067501029582579e-13"),
  @("13", "[[0.70769694]`n [0.63583678]]", "[[0.02064962]`n [0.09232268]]", [double]"9.657462472382125", [double]"0.609040853874799", [double]"370", [double]"1e-10", [double]"4.067501029582579e-13"),
  @("14", "[[0.76175856]`n [0.68229385]]", "[[0.02044505]`n [0.0946705 ]]", [double]"11.83951558400279", [double]"0.6073942367341977", [double]"292", [double]"1e-10", [double]"4.067501029582579e-13"),
  @("15", "[[0.81319885]`n [0.72189815]]", "[[0.019916  ]`n [0.09153467]]", [double]"13.88997841744493", [double]"0.6130269432859651", [double]"201", [double]"1e-10", [double]"4.067501029582579e-13"),
  @("16", "[[0.86076786]`n [0.76188949]]", "[[0.02046493]`n [0.1051492 ]]", [double]"15.82239359589182", [double]"0.5717476200783435", [double]"76", [double]"1e-10", [double]"4.067501029582579e-13"),
  @("17", "[[0.91426752]`n [0.81395877]]", "[[0.02075632]`n [0.10772341]]", [double]"18.00502137574836", [double]"0.5757820999538996", [double]"26", [double]"1e-10", [double]"4.067501029582579e-13"),
  @("18", "[[0.97319834]`n [0.85717933]]", "[[0.02503088]`n [0.12337471]]", [double]"20.33915870442998", [double]"0.8074960443910881", [double]"8", [double]"1e-10", [double]"4.067501029582579e-13")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 1).Value = [double]$rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
    $ws.Cells.Item($r, 5).Value = $rowVals[4]
    $ws.Cells.Item($r, 6).Value = $rowVals[5]
    $ws.Cells.Item($r, 7).Value = $rowVals[6]
    $ws.Cells.Item($r, 8).Value = $rowVals[7]
}

# Rows beyond the original 17 are brand-new cells and do not inherit the
# bold/bordered/centered header-column style used throughout column A;
# copy formatting only (keeps the numeric values already written) from an
# existing styled cell.
$ws.Range("A2").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Output "done"
